# The "<id>", "p128r_1" and "</id>" pieces of text live in three adjacent
# runs (Courier New / Arial / Courier New). The edit merges them into a
# single run - "<id>p128r_1</id>" - using the formatting of the first
# ("<id>") run, and removes the two runs that followed it.
#
# A Find & Replace across the combined range achieves exactly that: Word
# collapses the matched (multi-run) range into one run that carries the
# character formatting of the range's first character.

$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute(
    "<id>p128r_1</id>",   # find text
    $false,               # MatchCase
    $false,               # MatchWholeWord
    $false,               # MatchWildcards
    $false,               # MatchSoundsLike
    $false,               # MatchAllWordForms
    $true,                # Forward
    1,                    # Wrap (wdFindContinue)
    $false,               # Format
    "<id>p128r_1</id>",   # replacement text
    2                     # Replace (wdReplaceAll)
)

Write-Output "replaced: $found"
